$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (MAC address, serial number)
$macData = @{
    11 = @("3C:46:A1:2B:16:D0", "132373006295")
    12 = @("3C:46:A1:2A:FF:50", "132373005919")
    13 = @("3C:46:A1:2B:01:10", "132373005947")
    14 = @("3C:46:A1:2A:F8:C0", "132373005814")
    15 = @("3C:46:A1:2B:11:50", "132373006207")
    16 = @("3C:46:A1:2A:DB:30", "132373005341")
    17 = @("3C:46:A1:2B:37:50", "132373006815")
    18 = @("3C:46:A1:2B:27:80", "132373006562")
    19 = @("3C:46:A1:2A:D9:50", "132373005311")
    20 = @("3C:46:A1:2A:DB:F0", "132373005353")
    21 = @("3C:46:A1:2B:15:E0", "132373006280")
    22 = @("3C:46:A1:2B:29:80", "132373006594")
    23 = @("3C:46:A1:2B:36:A0", "132373006804")
    24 = @("3C:46:A1:2B:29:60", "132373006592")
    25 = @("3C:46:A1:2B:01:00", "132373005946")
    26 = @("3C:46:A1:2B:1E:40", "132373006414")
    48 = @("3C:46:A1:25:3C:10", "122379002745")
    49 = @("3C:46:A1:25:3A:F0", "122379002728")
    50 = @("3C:46:A1:25:3B:20", "122379002699")
    51 = @("3C:46:A1:25:3B:A0", "122379002737")
    52 = @("3C:46:A1:25:47:60", "122379002657")
    53 = @("3C:46:A1:25:47:B0", "122379002757")
    54 = @("3C:46:A1:25:78:30", "122379003294")
}

foreach ($row in $macData.Keys) {
    $vals = $macData[$row]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $vals[0]
    $cellC.Style = "Normal"

    # Column D holds a long numeric-looking serial number; force it to be
    # stored as text (matching the rest of the sheet) instead of a number,
    # then restore the default "Normal" style so formatting stays untouched.
    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $vals[1]
    $cellD.Style = "Normal"
}

# Rows 48-54: status column (F) changes from "libre" to "en uso"
foreach ($row in 48..54) {
    $ws.Cells.Item($row, 6).Value = "en uso"
}
